$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing hours for row 13 (Gabriela's Sunday time sheet)
$ws.Range("G13").Value = 2
$ws.Range("H13").Value = 4

# Reflect the final selected cell as seen in the saved file (H13)
$ws.Range("H13").Select()
